{"js": "// Add a \"blank line, then a line with today's date and your name\" and the\n// accompanying bio paragraph to the end of the document. The document's\n// final paragraph is already an empty paragraph (the \"blank line\"), so this\n// appends: <line break>date+name<line break>bio text, mirroring the\n// instructor's example entry that appears earlier in the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Target is the last paragraph in the document (the blank line left for\n// the student's entry).\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Run/font formatting used throughout the document body.\nconst rPr =\n  '<w:rPr><w:rFonts w:ascii=\"Calibri\" w:eastAsia=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Calibri\"/></w:rPr>';\n\nconst dateAndName = \"5/7/2024 Ceara Patterson\";\nconst bioIntro =\n  \"Hello all. I was born and raised here in Jacksonville. Although I\\u2019ve never lived anywhere else I do love to travel. I currently have plans to visit Vegas again, camp in the Blue Ridge, and am slowly hashing out plans for Ireland and Japan. I am an aspiring programmer in the Information Systems Tech bachelors program, interested in game development and front end web development. I enjoy camping/ hiking, video games,\";\nconst space = \" \";\nconst hobbies1 = \"anime, reading,\";\nconst hobbies2 = \" and gardening. \";\n\n// Small helper to XML-escape plain text before dropping it into <w:t>.\nfunction escapeXml(text) {\n  return text\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// <w:t> needs xml:space=\"preserve\" only when the text has leading/trailing\n// whitespace (matches how Word itself emits these runs).\nfunction textElement(text) {\n  const preserve = text !== text.trim() ? ' xml:space=\"preserve\"' : \"\";\n  return `<w:t${preserve}>${escapeXml(text)}</w:t>`;\n}\n\nconst runsXml =\n  `<w:r>${rPr}<w:br/>${textElement(dateAndName)}</w:r>` +\n  `<w:r>${rPr}<w:br/>${textElement(bioIntro)}</w:r>` +\n  `<w:r>${textElement(space)}</w:r>` +\n  `<w:r>${rPr}${textElement(hobbies1)}</w:r>` +\n  `<w:r>${rPr}${textElement(hobbies2)}</w:r>`;\n\n// Use an OOXML fragment (instead of plain insertText) so we can faithfully\n// reproduce the exact run/rPr structure (including w:cs, which the Word.js\n// Font object does not expose), matching what Word itself writes.\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  `<w:body><w:p>${runsXml}</w:p></w:body>` +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nlastParagraph.insertOoxml(ooxml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Add a \"blank line, then a line with today's date and your name\" and the\n# accompanying bio paragraph to the end of the document. The document's\n# final paragraph is already an empty paragraph (the \"blank line\"), so this\n# appends: <line break>date+name<line break>bio text, mirroring the\n# instructor's example entry that appears earlier in the document.\n\n$d = $word.ActiveDocument\n\n# Manual line break character (same as pressing Shift+Enter) -> <w:br/>.\n$lineBreak = [char]11\n\n$dateAndName = \"5/7/2024 Ceara Patterson\"\n$bioIntro = \"Hello all. I was born and raised here in Jacksonville. Although I\" + [char]0x2019 + \"ve never lived anywhere else I do love to travel. I currently have plans to visit Vegas again, camp in the Blue Ridge, and am slowly hashing out plans for Ireland and Japan. I am an aspiring programmer in the Information Systems Tech bachelors program, interested in game development and front end web development. I enjoy camping/ hiking, video games,\"\n$space = \" \"\n$hobbies1 = \"anime, reading,\"\n$hobbies2 = \" and gardening. \"\n\nfunction Set-CalibriFont($range) {\n    $range.Font.Name = \"Calibri\"\n    $range.Font.NameFarEast = \"Calibri\"\n    $range.Font.NameBi = \"Calibri\"\n}\n\n# Run 1: line break + date/name line.\n$r = $d.Paragraphs.Last.Range\n$r.Collapse(0)\n$r.InsertAfter($lineBreak + $dateAndName)\nSet-CalibriFont $r\n\n# Run 2: line break + bio intro (through \"...video games,\").\n$r = $d.Paragraphs.Last.Range\n$r.Collapse(0)\n$r.InsertAfter($lineBreak + $bioIntro)\nSet-CalibriFont $r\n\n# Run 3: a lone space, left with no explicit run formatting (matches the\n# source document, where this run carries no rPr).\n$r = $d.Paragraphs.Last.Range\n$r.Collapse(0)\n$r.InsertAfter($space)\n\n# Run 4: \"anime, reading,\"\n$r = $d.Paragraphs.Last.Range\n$r.Collapse(0)\n$r.InsertAfter($hobbies1)\nSet-CalibriFont $r\n\n# Run 5: \" and gardening. \"\n$r = $d.Paragraphs.Last.Range\n$r.Collapse(0)\n$r.InsertAfter($hobbies2)\nSet-CalibriFont $r\n"}
